# Applies the "16th May Refresh" master-data update:
#  1. Renames three existing template codes (the otp-* codes become
#     ida-auth-otp-*-template, reusing their existing descriptions).
#  2. Appends 12 new rows (125-136) for the new "consent" and
#     "auth-otp-*-template" master data entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the three existing codes (column A only) ------------------
# otp-email-content  -> ida-auth-otp-email-content-template   (rows 5,11,17)
# otp-email-subject  -> ida-auth-otp-email-subject-template   (rows 6,12,18)
# otp-sms            -> ida-auth-otp-sms-template              (rows 7,13,19)
$ws.Range("A5").Value  = "ida-auth-otp-email-content-template"
$ws.Range("A11").Value = "ida-auth-otp-email-content-template"
$ws.Range("A17").Value = "ida-auth-otp-email-content-template"

$ws.Range("A6").Value  = "ida-auth-otp-email-subject-template"
$ws.Range("A12").Value = "ida-auth-otp-email-subject-template"
$ws.Range("A18").Value = "ida-auth-otp-email-subject-template"

$ws.Range("A7").Value  = "ida-auth-otp-sms-template"
$ws.Range("A13").Value = "ida-auth-otp-sms-template"
$ws.Range("A19").Value = "ida-auth-otp-sms-template"

# --- 2. Append the new rows 125-136 --------------------------------------
$newRows = @(
    @("consent", "Consent", "eng"),
    @("consent", "موافقة", "ara"),
    @("consent", "Consentement", "fra"),
    @("auth-otp-email-subject-template", "Auth OTP Email Subject Template", "eng"),
    @("auth-otp-email-subject-template", "مصادقة OTP قالب موضوع", "ara"),
    @("auth-otp-email-subject-template", "Modèle dobjet de-mail Auth OTP", "fra"),
    @("auth-otp-email-content-template", "Auth OTP Email Content Template", "eng"),
    @("auth-otp-email-content-template", "مصادقة OTP قالب محتوى", "ara"),
    @("auth-otp-email-content-template", "Auth OTP Email ContentTemplate", "fra"),
    @("auth-otp-sms-template", "Auth OTP SMS Template", "eng"),
    @("auth-otp-sms-template", "مصادقة قالب رسالة OTP", "ara"),
    @("auth-otp-sms-template", "Modèle SMS OTP Auth", "fra")
)

$row = 125
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $row++
}

# --- 3. Refresh the "new row" selection marker, like Excel does ----------
$ws.Range("A137:XFD1048576").Select()
